# Apply 2025-10-28 data updates (2025 column = 'L', plus 2021 column = 'H' for a couple totals)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5540
$ws.Range("L3").Value = 6036
$ws.Range("H4").Value = 1765
$ws.Range("L4").Value = 1482
$ws.Range("L5").Value = 361
$ws.Range("L6").Value = 4956
$ws.Range("H7").Value = 26081
$ws.Range("L7").Value = 18375

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 428
$ws.Range("L7").Value = 1214

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 166
$ws.Range("L7").Value = 408

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L4").Value = 56
$ws.Range("L7").Value = 845

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 245
$ws.Range("L5").Value = 21
$ws.Range("L6").Value = 191
$ws.Range("L7").Value = 704

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 113
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 354

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 130
$ws.Range("L7").Value = 319

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 144
$ws.Range("L7").Value = 594
$ws.Range("L8").Value = 1214
$ws.Range("L12").Value = 43
$ws.Range("L15").Value = 149
$ws.Range("L16").Value = 40
$ws.Range("L17").Value = 35
$ws.Range("L18").Value = 124
$ws.Range("L19").Value = 501
$ws.Range("L20").Value = 453
$ws.Range("L22").Value = 55
$ws.Range("L23").Value = 200
$ws.Range("L29").Value = 1036
$ws.Range("L30").Value = 81
$ws.Range("L33").Value = 845
$ws.Range("L34").Value = 107
$ws.Range("L37").Value = 704
$ws.Range("L40").Value = 52
$ws.Range("L42").Value = 597
$ws.Range("L43").Value = 132
$ws.Range("L44").Value = 124
$ws.Range("L46").Value = 41
$ws.Range("L48").Value = 240
$ws.Range("L50").Value = 90
$ws.Range("L53").Value = 201
$ws.Range("L54").Value = 403
$ws.Range("L55").Value = 194
$ws.Range("H63").Value = 316
$ws.Range("L65").Value = 354
$ws.Range("L67").Value = 633
$ws.Range("L73").Value = 147
$ws.Range("L78").Value = 233
$ws.Range("L79").Value = 500
$ws.Range("L83").Value = 408
$ws.Range("L87").Value = 55
$ws.Range("L91").Value = 243
$ws.Range("L94").Value = 225
$ws.Range("L96").Value = 207
$ws.Range("L97").Value = 152
$ws.Range("L99").Value = 319
$ws.Range("H101").Value = 26081
$ws.Range("L101").Value = 18375

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 245
$ws.Range("L7").Value = 633

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 103
$ws.Range("L6").Value = 193
$ws.Range("L7").Value = 403

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 307
$ws.Range("L6").Value = 260
$ws.Range("L7").Value = 1036

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 48
$ws.Range("L7").Value = 240

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 176
$ws.Range("L3").Value = 155
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 501

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L3").Value = 42
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 162
$ws.Range("L7").Value = 597

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 77
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 194

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 79
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 82
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 500

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 153
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 453

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 594

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 52
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 56
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L6").Value = 71
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 376
$ws.Range("L4").Value = 52

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 40
